# Auto-generated edit script applying numeric updates from the commit diff
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2038.1364
$ws.Range("I28").Value = 1717.7858
$ws.Range("J28").Value = 2598.75
$ws.Range("K28").Value = 1717.7858
$ws.Range("L28").Value = 2598.75
$ws.Range("M28").Value = -1232.7858
$ws.Range("N28").Value = -3568.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 601.8214
$ws.Range("I33").Value = 349.26086
$ws.Range("J33").Value = 1763.6
$ws.Range("K33").Value = 349.26086
$ws.Range("L33").Value = 1763.6
$ws.Range("M33").Value = -120.26086
$ws.Range("N33").Value = -2221.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3063.0435
$ws.Range("I64").Value = 2788.4614
$ws.Range("K64").Value = 2788.4614
$ws.Range("M64").Value = -2540.4614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3063.0435
$ws.Range("I67").Value = 2788.4614
$ws.Range("K67").Value = 2788.4614
$ws.Range("M67").Value = -1930.4614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5467.25
$ws.Range("J86").Value = 6817.3335
$ws.Range("L86").Value = 6817.3335
$ws.Range("N86").Value = -9063.333500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 5467.25
$ws.Range("J89").Value = 6817.3335
$ws.Range("L89").Value = 34086.6675
$ws.Range("N89").Value = -45318.6675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3000
$ws.Range("J116").Value = 3000
$ws.Range("L116").Value = 3000
$ws.Range("N116").Value = -9884

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 583.73334
$ws.Range("I121").Value = 263
$ws.Range("J121").Value = 633.0769
$ws.Range("K121").Value = 789
$ws.Range("L121").Value = 1899.2307
$ws.Range("M121").Value = 958
$ws.Range("N121").Value = -5393.2307

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4991.2905
$ws.Range("I132").Value = 1737.36
$ws.Range("J132").Value = 18549.334
$ws.Range("K132").Value = 5212.08
$ws.Range("L132").Value = 55648.00199999999
$ws.Range("M132").Value = -2682.08
$ws.Range("N132").Value = -60708.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12951.333
$ws.Range("I32").Value = 7876.109
$ws.Range("J32").Value = 23101.783
$ws.Range("K32").Value = 7876.109
$ws.Range("L32").Value = 23101.783
$ws.Range("M32").Value = -7589.109
$ws.Range("N32").Value = -23675.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M61").Value = -1005.5
$ws.Range("N61").ClearContents()
$ws.Range("H61").Value = 1217.5
$ws.Range("I61").Value = 1217.5
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1217.5
$ws.Range("L61").Value = 0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3898
$ws.Range("I63").Value = 2995
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 2995
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -2309
$ws.Range("N63").Value = -5872

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3898
$ws.Range("I66").Value = 2995
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 14975
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -11543
$ws.Range("N66").Value = -29364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M136").Value = -1102.5
$ws.Range("N136").ClearContents()
$ws.Range("H136").Value = 1217.5
$ws.Range("I136").Value = 1217.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 3652.5
$ws.Range("L136").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20001822
$ws.Range("I20").Value = 30304700
$ws.Range("J20").Value = 2120.9412
$ws.Range("K20").Value = 30304700
$ws.Range("L20").Value = 2120.9412
$ws.Range("M20").Value = -30304453
$ws.Range("N20").Value = -2614.9412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1981.174
$ws.Range("I86").Value = 1851.2941
$ws.Range("J86").Value = 2349.1667
$ws.Range("K86").Value = 1851.2941
$ws.Range("L86").Value = 2349.1667
$ws.Range("M86").Value = -728.2941000000001
$ws.Range("N86").Value = -4595.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1981.174
$ws.Range("I89").Value = 1851.2941
$ws.Range("J89").Value = 2349.1667
$ws.Range("K89").Value = 9256.470499999999
$ws.Range("L89").Value = 11745.8335
$ws.Range("M89").Value = -3640.470499999999
$ws.Range("N89").Value = -22977.8335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1796.75
$ws.Range("I94").Value = 1366.0714
$ws.Range("K94").Value = 1366.0714
$ws.Range("M94").Value = -915.0714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2388.3965
$ws.Range("I105").Value = 2385.3962
$ws.Range("K105").Value = 2385.3962
$ws.Range("M105").Value = -638.3962000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750
$ws.Range("J16").Value = 900
$ws.Range("L16").Value = 900
$ws.Range("N16").Value = -1474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12499.75
$ws.Range("J50").Value = 12499.75
$ws.Range("L50").Value = 12499.75
$ws.Range("N50").Value = -13749.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 12999.5
$ws.Range("J51").Value = 12999.5
$ws.Range("L51").Value = 12999.5
$ws.Range("N51").Value = -14471.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 16909.092
$ws.Range("J59").Value = 16909.092
$ws.Range("L59").Value = 16909.092
$ws.Range("N59").Value = -19199.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 6459.8
$ws.Range("I60").Value = 700
$ws.Range("J60").Value = 7899.75
$ws.Range("K60").Value = 700
$ws.Range("L60").Value = 7899.75
$ws.Range("M60").Value = -189
$ws.Range("N60").Value = -8921.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 12999.5
$ws.Range("J61").Value = 12999.5
$ws.Range("L61").Value = 12999.5
$ws.Range("N61").Value = -13695.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2522.7273
$ws.Range("I62").Value = 2512.5
$ws.Range("J62").Value = 2550
$ws.Range("K62").Value = 2512.5
$ws.Range("L62").Value = 2550
$ws.Range("M62").Value = -1888.5
$ws.Range("N62").Value = -3798

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2522.7273
$ws.Range("I65").Value = 2512.5
$ws.Range("J65").Value = 2550
$ws.Range("K65").Value = 12562.5
$ws.Range("L65").Value = 12750
$ws.Range("M65").Value = -9442.5
$ws.Range("N65").Value = -18990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 750
$ws.Range("J113").Value = 900
$ws.Range("L113").Value = 900
$ws.Range("N113").Value = -5240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 700759.3
$ws.Range("I132").Value = 1756.3077
$ws.Range("J132").Value = 2647982
$ws.Range("K132").Value = 5268.9231
$ws.Range("L132").Value = 7943946
$ws.Range("M132").Value = -2738.9231
$ws.Range("N132").Value = -7949006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 8000
$ws.Range("I38").Value = 5000
$ws.Range("J38").Value = 11000
$ws.Range("K38").Value = 5000
$ws.Range("L38").Value = 11000
$ws.Range("M38").Value = -4537
$ws.Range("N38").Value = -11926

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 10714.286
$ws.Range("J40").Value = 10714.286
$ws.Range("L40").Value = 10714.286
$ws.Range("N40").Value = -11016.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6212.5
$ws.Range("I70").Value = 5540
$ws.Range("K70").Value = 5540
$ws.Range("M70").Value = -5270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 6212.5
$ws.Range("I73").Value = 5540
$ws.Range("K73").Value = 5540
$ws.Range("M73").Value = -4604

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2588.889
$ws.Range("I80").Value = 2325
$ws.Range("J80").Value = 2800
$ws.Range("K80").Value = 2325
$ws.Range("L80").Value = 2800
$ws.Range("M80").Value = -1327
$ws.Range("N80").Value = -4796

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2588.889
$ws.Range("I83").Value = 2325
$ws.Range("J83").Value = 2800
$ws.Range("K83").Value = 11625
$ws.Range("L83").Value = 14000
$ws.Range("M83").Value = -6633
$ws.Range("N83").Value = -23984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 35128.45
$ws.Range("I22").Value = 111564.445
$ws.Range("J22").Value = 732.25
$ws.Range("K22").Value = 111564.445
$ws.Range("L22").Value = 732.25
$ws.Range("M22").Value = -111269.445
$ws.Range("N22").Value = -1322.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 35128.45
$ws.Range("I27").Value = 111564.445
$ws.Range("J27").Value = 732.25
$ws.Range("K27").Value = 111564.445
$ws.Range("L27").Value = 732.25
$ws.Range("M27").Value = -111457.445
$ws.Range("N27").Value = -946.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2150.125
$ws.Range("I122").Value = 1321.2858
$ws.Range("K122").Value = 3963.8574
$ws.Range("M122").Value = -1513.8574

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 341.73172
$ws.Range("I113").Value = 290.2
$ws.Range("K113").Value = 870.5999999999999
$ws.Range("M113").Value = 1299.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3582.9583
$ws.Range("I132").Value = 3630.875
$ws.Range("J132").Value = 3487.125
$ws.Range("K132").Value = 10892.625
$ws.Range("L132").Value = 10461.375
$ws.Range("M132").Value = -8362.625
$ws.Range("N132").Value = -15521.375

Write-Host "Applied all profit-sheet updates"
